# Add 2022-Q4 data
# -----------------------------------------------------------------------
# Before: workbook has 2 sheets -> "总计" (summary) and "2021-Q4" (fund
# holdings for that quarter).
# After:  workbook has 3 sheets -> "总计", "2022-Q4" (new fund holdings),
# "2021-Q4" (unchanged fund holdings, moved to the last position).
# The summary sheet gets a new row for 2022-Q4 inserted above the
# existing 2021-Q4 summary row.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)      # "总计"
$oldQ  = $wb.Worksheets.Item(2)      # currently "2021-Q4", holds the fund data

# -----------------------------------------------------------------------
# 1. Create a brand-new sheet right after the existing quarter sheet, and
#    move the *old* 2021-Q4 fund data into it (so that the original
#    "2021-Q4" sheet content/ids are preserved, and the in-place sheet
#    can be repurposed for the new 2022-Q4 numbers).
# -----------------------------------------------------------------------
$oldQ.Name = "2021-Q4 (old)"
$newQ = $wb.Worksheets.Add($null, $oldQ)
$newQ.Name = "2021-Q4"

# Copy header row (B1:H1) formatting + values from the old sheet.
$oldQ.Range("B1:H1").Copy()
$newQ.Range("B1:H1").PasteSpecial(-4122)
$newQ.Range("B1").Value2 = $oldQ.Range("B1").Value2
$newQ.Range("C1").Value2 = $oldQ.Range("C1").Value2
$newQ.Range("D1").Value2 = $oldQ.Range("D1").Value2
$newQ.Range("E1").Value2 = $oldQ.Range("E1").Value2
$newQ.Range("F1").Value2 = $oldQ.Range("F1").Value2
$newQ.Range("G1").Value2 = $oldQ.Range("G1").Value2
$newQ.Range("H1").Value2 = $oldQ.Range("H1").Value2

# Copy row 2 (fund holding values) formatting, then the values themselves.
$oldQ.Range("A2:H2").Copy()
$newQ.Range("A2:H2").PasteSpecial(-4122)

$newQ.Range("A2").Value2 = 0
$newQ.Range("B2:G2").NumberFormat = "@"
$newQ.Range("B2").Value2 = $oldQ.Range("B2").Value2
$newQ.Range("C2").Value2 = $oldQ.Range("C2").Value2
$newQ.Range("D2").Value2 = $oldQ.Range("D2").Value2
$newQ.Range("E2").Value2 = $oldQ.Range("E2").Value2
$newQ.Range("F2").Value2 = $oldQ.Range("F2").Value2
$newQ.Range("G2").Value2 = $oldQ.Range("G2").Value2
$newQ.Range("H2").Value2 = 4

# -----------------------------------------------------------------------
# 2. Re-purpose the in-place sheet (still rId/sheetId of the original
#    "2021-Q4" sheet) as the new "2022-Q4" sheet, with the new fund data.
#    Row 1 headers are identical, so they stay untouched.
# -----------------------------------------------------------------------
$oldQ.Name = "2022-Q4"

$oldQ.Range("A2").Value2 = 0
$oldQ.Range("B2:G2").NumberFormat = "@"
$oldQ.Range("B2").Value2 = "012977"
$oldQ.Range("C2").Value2 = "瑞达鑫红量化6个月持有混合A"
$oldQ.Range("D2").Value2 = "0.35"
$oldQ.Range("E2").Value2 = "94.66"
$oldQ.Range("F2").Value2 = "4.94"
$oldQ.Range("G2").Value2 = "0.0173"
$oldQ.Range("H2").Value2 = 6

# New row 3 - copy formatting from row 2, then fill in the 012978 entry.
$oldQ.Range("A2:H2").Copy()
$oldQ.Range("A3:H3").PasteSpecial(-4122)

$oldQ.Range("A3").Value2 = 1
$oldQ.Range("B3:G3").NumberFormat = "@"
$oldQ.Range("B3").Value2 = "012978"
$oldQ.Range("C3").Value2 = "瑞达鑫红量化6个月持有混合C"
$oldQ.Range("D3").Value2 = "0.09"
$oldQ.Range("E3").Value2 = "94.66"
$oldQ.Range("F3").Value2 = "4.94"
$oldQ.Range("G3").Value2 = "0.0044"
$oldQ.Range("H3").Value2 = 6

# -----------------------------------------------------------------------
# 3. Update the "总计" sheet: row 2 becomes the 2022-Q4 summary, and a new
#    row 3 is appended with the 2021-Q4 summary (same values that used to
#    live in row 2).
# -----------------------------------------------------------------------
$total.Range("A2:D2").Copy()
$total.Range("A3:D3").PasteSpecial(-4122)

$total.Range("A3").Value2 = 1
$total.Range("B3").Value2 = "2021-Q4"
$total.Range("C3").Value2 = 1
$total.Range("D3").Value2 = 0.85

$total.Range("A2").Value2 = 0
$total.Range("B2").Value2 = "2022-Q4"
$total.Range("C2").Value2 = 2
$total.Range("D2").Value2 = 0.02
